# "Generate Report for Handback"
#
# For each locale sheet (zh-cn, de-de) the handback engine has now produced a
# target file and a handback file for both tracked source documents, so the
# report gains:
#   - "Latest Target File"   (col I) -> hyperlink to the source .md file
#   - "Latest Handback File" (col J) -> the locale-specific handback .xlf name
#   - "Latest Handback DateTime" (col K) -> the handback timestamp
# The Overview sheet's "Latest HO Xliff Generate Date" status text also flips
# from "In Translation" to "Handed back: in sync with en-US" for both rows,
# and a handful of columns widen to comfortably fit the long file names.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: status column (E & F) now reports the handback sync state.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.1
$overview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------------
# zh-cn sheet: populate target/handback file + datetime for both rows.
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("J2").Value = "83f96737-3f85-4cda-8076-2bed71007eda.a2ae9151179e7e297ff7b3011d9deeeaf1a54e1c.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-09-01 04:28:13"

$zhcn.Range("J3").Value = "e013ab18-ac30-45bf-a6cb-a8be1a7758a9.25ba904a8566c274f1e375dc6a971efe4ee60812.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-09-01 04:28:13"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e761382c1f381599dfe2a0c8405be2c9cd202e3a/e2e/83f96737-3f85-4cda-8076-2bed71007eda.md", "", "", "83f96737-3f85-4cda-8076-2bed71007eda.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e761382c1f381599dfe2a0c8405be2c9cd202e3a/e2e/e013ab18-ac30-45bf-a6cb-a8be1a7758a9.md", "", "", "e013ab18-ac30-45bf-a6cb-a8be1a7758a9.md") | Out-Null

$zhcn.Columns.Item(3).ColumnWidth = 29.1
$zhcn.Columns.Item(9).ColumnWidth = 39.17
$zhcn.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet: same shape, different locale file names/timestamp.
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("J2").Value = "83f96737-3f85-4cda-8076-2bed71007eda.a2ae9151179e7e297ff7b3011d9deeeaf1a54e1c.de-de.xlf"
$dede.Range("K2").Value = "2016-09-01 04:28:21"

$dede.Range("J3").Value = "e013ab18-ac30-45bf-a6cb-a8be1a7758a9.25ba904a8566c274f1e375dc6a971efe4ee60812.de-de.xlf"
$dede.Range("K3").Value = "2016-09-01 04:28:21"

$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e761382c1f381599dfe2a0c8405be2c9cd202e3a/e2e/83f96737-3f85-4cda-8076-2bed71007eda.md", "", "", "83f96737-3f85-4cda-8076-2bed71007eda.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e761382c1f381599dfe2a0c8405be2c9cd202e3a/e2e/e013ab18-ac30-45bf-a6cb-a8be1a7758a9.md", "", "", "e013ab18-ac30-45bf-a6cb-a8be1a7758a9.md") | Out-Null

$dede.Columns.Item(3).ColumnWidth = 29.1
$dede.Columns.Item(9).ColumnWidth = 39.17
$dede.Columns.Item(10).ColumnWidth = 39.17
